$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-11-22 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-23 Sunday", 2) | Out-Null

# Update the math table cells (positional replace, row by row, left to right)
$t = $d.Tables(1)
$t.Cell(1, 1).Range.Text = "8+8="
$t.Cell(1, 2).Range.Text = "93-24="
$t.Cell(1, 3).Range.Text = "86+9="
$t.Cell(1, 4).Range.Text = "98-9="
$t.Cell(1, 5).Range.Text = "4+27="
$t.Cell(2, 1).Range.Text = "65+28="
$t.Cell(2, 2).Range.Text = "43-15="
$t.Cell(2, 3).Range.Text = "47+44="
$t.Cell(2, 4).Range.Text = "61-4="
$t.Cell(2, 5).Range.Text = "19+4="
$t.Cell(3, 1).Range.Text = "5+88="
$t.Cell(3, 2).Range.Text = "6+35="
$t.Cell(3, 3).Range.Text = "59+6="
$t.Cell(3, 4).Range.Text = "75-58="
$t.Cell(3, 5).Range.Text = "85-7="
$t.Cell(4, 1).Range.Text = "40-23="
$t.Cell(4, 2).Range.Text = "48+37="
$t.Cell(4, 3).Range.Text = "91-56="
$t.Cell(4, 4).Range.Text = "47+28="
$t.Cell(4, 5).Range.Text = "49+36="
$t.Cell(5, 1).Range.Text = "6+78="
$t.Cell(5, 2).Range.Text = "62+9="
$t.Cell(5, 3).Range.Text = "14+28="
$t.Cell(5, 4).Range.Text = "70-23="
$t.Cell(5, 5).Range.Text = "9+77="
$t.Cell(6, 1).Range.Text = "28+4="
$t.Cell(6, 2).Range.Text = "37+5="
$t.Cell(6, 3).Range.Text = "90-15="
$t.Cell(6, 4).Range.Text = "75-66="
$t.Cell(6, 5).Range.Text = "39+24="
$t.Cell(7, 1).Range.Text = "69+5="
$t.Cell(7, 2).Range.Text = "31-8="
$t.Cell(7, 3).Range.Text = "46-37="
$t.Cell(7, 4).Range.Text = "34+58="
$t.Cell(7, 5).Range.Text = "94-17="
$t.Cell(8, 1).Range.Text = "29+49="
$t.Cell(8, 2).Range.Text = "15+77="
$t.Cell(8, 3).Range.Text = "17+35="
$t.Cell(8, 4).Range.Text = "94-6="
$t.Cell(8, 5).Range.Text = "10-2="
$t.Cell(9, 1).Range.Text = "56-17="
$t.Cell(9, 2).Range.Text = "41-14="
$t.Cell(9, 3).Range.Text = "35+36="
$t.Cell(9, 4).Range.Text = "36-29="
$t.Cell(9, 5).Range.Text = "5+36="
$t.Cell(10, 1).Range.Text = "53+18="
$t.Cell(10, 2).Range.Text = "33+9="
$t.Cell(10, 3).Range.Text = "46+37="
$t.Cell(10, 4).Range.Text = "68-59="
$t.Cell(10, 5).Range.Text = "9+74="
$t.Cell(11, 1).Range.Text = "61-44="
$t.Cell(11, 3).Range.Text = "32+39="
$t.Cell(11, 4).Range.Text = "26-18="
$t.Cell(11, 5).Range.Text = "55+16="
$t.Cell(12, 1).Range.Text = "8+23="
$t.Cell(12, 2).Range.Text = "71-59="
$t.Cell(12, 3).Range.Text = "57+26="
$t.Cell(12, 4).Range.Text = "84-79="
$t.Cell(12, 5).Range.Text = "86+7="
$t.Cell(13, 1).Range.Text = "56+25="
$t.Cell(13, 2).Range.Text = "29+33="
$t.Cell(13, 3).Range.Text = "95-88="
$t.Cell(13, 4).Range.Text = "18+38="
$t.Cell(13, 5).Range.Text = "63-47="
$t.Cell(14, 1).Range.Text = "81-34="
$t.Cell(14, 2).Range.Text = "63-54="
$t.Cell(14, 3).Range.Text = "84-79="
$t.Cell(14, 4).Range.Text = "66+5="
$t.Cell(14, 5).Range.Text = "7+74="
$t.Cell(15, 1).Range.Text = "34-16="
$t.Cell(15, 2).Range.Text = "58+19="
$t.Cell(15, 3).Range.Text = "7+37="
$t.Cell(15, 4).Range.Text = "66+6="
$t.Cell(15, 5).Range.Text = "17+35="
$t.Cell(16, 1).Range.Text = "25-17="
$t.Cell(16, 2).Range.Text = "91-17="
$t.Cell(16, 3).Range.Text = "47+39="
$t.Cell(16, 4).Range.Text = "9+49="
$t.Cell(16, 5).Range.Text = "13+48="
$t.Cell(17, 1).Range.Text = "47+17="
$t.Cell(17, 2).Range.Text = "91-7="
$t.Cell(17, 3).Range.Text = "70-29="
$t.Cell(17, 4).Range.Text = "66+27="
$t.Cell(17, 5).Range.Text = "9+66="
$t.Cell(18, 1).Range.Text = "34+28="
$t.Cell(18, 2).Range.Text = "6+35="
$t.Cell(18, 3).Range.Text = "84-5="
$t.Cell(18, 4).Range.Text = "82-5="
$t.Cell(18, 5).Range.Text = "58+24="
$t.Cell(19, 1).Range.Text = "42-39="
$t.Cell(19, 2).Range.Text = "53-9="
$t.Cell(19, 3).Range.Text = "57+9="
$t.Cell(19, 4).Range.Text = "42-34="
$t.Cell(19, 5).Range.Text = "49+37="
$t.Cell(20, 1).Range.Text = "27+37="
$t.Cell(20, 2).Range.Text = "80-56="
$t.Cell(20, 3).Range.Text = "95-7="
$t.Cell(20, 4).Range.Text = "48+44="
$t.Cell(20, 5).Range.Text = "49+2="
